$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC2").Value = 19
$ws.Range("AF2").Value = 19

# Row 3
$ws.Range("D3").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC3").Value = 19
$ws.Range("AF3").Value = 19

# Row 4
$ws.Range("D4").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC4").Value = 18
$ws.Range("AF4").Value = 18

# Row 5
$ws.Range("D5").Value = "2024-07-19T12:51:00.000Z"

# Row 6
$ws.Range("D6").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC6").Value = 19
$ws.Range("AF6").Value = 19

# Row 7
$ws.Range("D7").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC7").Value = 19
$ws.Range("AF7").Value = 19

# Row 8
$ws.Range("D8").Value = "2024-07-19T12:51:00.000Z"

# Row 9
$ws.Range("D9").Value = "2024-07-19T12:51:00.000Z"
$ws.Range("AC9").Value = 18
$ws.Range("AF9").Value = 19

# Row 10
$ws.Range("D10").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC10").Value = 18
$ws.Range("AF10").Value = 18

# Row 11
$ws.Range("D11").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("S11").Value = 2
$ws.Range("AC11").Value = 19
$ws.Range("AF11").Value = 21

# Row 12
$ws.Range("D12").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC12").Value = 17
$ws.Range("AF12").Value = 19

# Row 13
$ws.Range("D13").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC13").Value = 14
$ws.Range("AF13").Value = 14

# Row 14
$ws.Range("D14").Value = "2024-07-19T12:52:00.000Z"

# Row 15
$ws.Range("D15").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC15").Value = 16
$ws.Range("AF15").Value = 17

# Row 16
$ws.Range("D16").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC16").Value = 13
$ws.Range("AF16").Value = 18

# Row 17
$ws.Range("D17").Value = "2024-07-19T12:52:00.000Z"

# Row 18
$ws.Range("D18").Value = "2024-07-19T12:52:00.000Z"
$ws.Range("AC18").Value = 19
$ws.Range("AF18").Value = 19

# Row 19
$ws.Range("D19").Value = "2024-07-19T12:52:00.000Z"

# Row 20
$ws.Range("D20").Value = "2024-07-19T12:52:00.000Z"
